$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 (XKS) status goes from "Incomplete" (red) to "In Progress" (yellow).
# Capture F13's current "In Progress" (yellow) formatting before it changes below.
$ws.Range("F13").Copy() | Out-Null
$ws.Range("F23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F23").Value = "In Progress"

# Row 13 (IPStats) status goes from "In Progress" (yellow) to the new
# "Complete 0.3.1.6b" (green), matching the other "Complete" cells like F16.
$ws.Range("F16").Copy() | Out-Null
$ws.Range("F13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F13").Value = "Complete 0.3.1.6b"

$excel.CutCopyMode = 0

# Move the active selection to F16.
$ws.Range("F16").Select() | Out-Null
